# Fix example placeholder formatting: strip the {{ }} mustache wrapper
# around the template tokens, leaving the bare token names.

$d = $word.ActiveDocument

$tokens = @(
    "CONTRACT_DATE",
    "CLIENT_NAME",
    "CLIENT_COMPANY",
    "CLIENT_EMAIL",
    "CLIENT_PHONE",
    "PROJECT_NAME",
    "CONTRACT_AMOUNT",
    "PROJECT_DEADLINE",
    "PAYMENT_TERMS"
)

foreach ($token in $tokens) {
    $find = "{{" + $token + "}}"
    $rng = $d.Content
    # MatchCase=True, MatchWholeWord=False, MatchWildcards=False,
    # MatchSoundsLike=False, MatchAllWordForms=False, Forward=True,
    # Wrap=wdFindContinue(1), Format=False, ReplaceWith=$token,
    # Replace=wdReplaceAll(2) — replaces every occurrence in the doc.
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $token, 2)
}

# Even out the template table's column widths (Width is in points;
# gridCol stores twips = points * 20).
$tbl = $d.Tables.Item(1)
$tbl.Columns.Item(1).Width = 2477 / 20
$tbl.Columns.Item(2).Width = 2478 / 20
$tbl.Columns.Item(3).Width = 2478 / 20
$tbl.Columns.Item(4).Width = 2478 / 20
